$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 14
$ws.Range("B14").Value = 8
$ws.Range("C14").Value = 'Javascript action listener should detect use of the submit button'
$ws.Range("D14").Value = 'Submit button and javascript actionlistener need to be present'
$ws.Range("E14").Value = 'console.log("HELLO :)");'
$ws.Range("F14").Value = 'After adding function, check submit button produces test data in console'
$ws.Range("G14").Value = 'hello in console'
$ws.Range("H14").Value = 'Console window epxlaining wrong syntax'
$ws.Range("I14").Value = 'fail'
$ws.Range("K14").Value = 'https://imgur.com/O9AQ81t'
$ws.Rows.Item(14).RowHeight = 60

# Row 15
$ws.Range("B15").Value = 9
$ws.Range("C15").Value = 'Javascript action listener should detect user click input'
$ws.Range("D15").Value = 'javascript actionlistener needs to be present'
$ws.Range("E15").Value = 'console.log("HELLO :)");'
$ws.Range("F15").Value = 'After adding function, check submit button produces test data in console'
$ws.Range("G15").Value = 'hello in console'
$ws.Range("H15").Value = 'Shows hello, but no functionallity.'
$ws.Range("I15").Value = 'fail'
$ws.Range("J15").Value = 'Text loads when script loads. Believe this could be due to loading in the head, instead of body but also could be due to the way I have made the function'
$ws.Range("K15").Value = 'https://imgur.com/GCRzhH9'
$ws.Rows.Item(15).RowHeight = 45

# Row 16
$ws.Range("B16").Value = 10
$ws.Range("C16").Value = 'Javascript action listener should detect user click input'
$ws.Range("D16").Value = 'javascript actionlistener needs to be present'
$ws.Range("E16").Value = 'console.log("HELLO :)");'
$ws.Range("F16").Value = 'While using an eventlistener function found online click on the screen'
$ws.Range("G16").Value = 'hello in console'
$ws.Range("H16").Value = 'Produces hello in console when clicked anywhere on document'
$ws.Range("I16").Value = 'pass'
$ws.Range("J16").Value = 'I think I was trying to incorporate one of these =>, since it was what I used last time. Also code is from w3schools'
$ws.Range("K16").Value = 'https://imgur.com/2cPukwi'
$ws.Rows.Item(16).RowHeight = 45

# Row 17
$ws.Range("B17").Value = 10
$ws.Range("C17").Value = 'Javascript action listener should detect use of the submit button'
$ws.Range("D17").Value = 'javascript actionlistener needs to be present'
$ws.Range("E17").Value = 'console.log("HELLO :)");'
$ws.Range("F17").Value = 'After editing old  function, check submit button produces test data in console'
$ws.Range("G17").Value = 'hello in console'
$ws.Range("H17").Value = 'Produces hello in console when submit button is clicked'
$ws.Range("I17").Value = 'pass'
$ws.Range("K17").Value = 'https://imgur.com/QVmL1EK'
$ws.Rows.Item(17).RowHeight = 45

# Row 18
$ws.Range("B18").Value = 10
$ws.Range("C18").Value = 'All fields detect user input (click)'
$ws.Range("D18").Value = 'actionlisteners + html elements'
$ws.Range("E18").Value = 'console.log("HELLO :)");'
$ws.Range("F18").Value = 'Add all eventlisteners to all inputs'
$ws.Range("G18").Value = 'hello in console for each input clicked on'
$ws.Range("H18").Value = 'All fields produce hello when clicked '
$ws.Range("I18").Value = 'pass'
$ws.Range("J18").Value = 'Used an array and looped through to add event listeners to each. Think this might need changing for code readablity, since userElement[1] isn''t very descriptive'
$ws.Range("K18").Value = 'https://imgur.com/DJJxcTA'
$ws.Rows.Item(18).RowHeight = 45

# Hyperlinks for K14:K18
$ws.Hyperlinks.Add($ws.Range("K14"), 'https://imgur.com/O9AQ81t')
$ws.Hyperlinks.Add($ws.Range("K15"), 'https://imgur.com/GCRzhH9')
$ws.Hyperlinks.Add($ws.Range("K16"), 'https://imgur.com/2cPukwi')
$ws.Hyperlinks.Add($ws.Range("K17"), 'https://imgur.com/QVmL1EK')
$ws.Hyperlinks.Add($ws.Range("K18"), 'https://imgur.com/DJJxcTA')
$ws.Range("K13").Copy()
$ws.Range("K14:K18").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Update selection
$ws.Range("I16").Select()
